$wb = $excel.ActiveWorkbook

# --- Rename sheets (scores -> score, parameters -> param, indicators -> ind) ---
$wb.Worksheets.Item("scores").Name = "score"
$wb.Worksheets.Item("parameters").Name = "param"
$wb.Worksheets.Item("indicators").Name = "ind"

# --- Fix defined names: #NAME? -> #name? (lowercase error token) ---
$wb.Names.Item("dbl_BenchmarkMax").RefersTo = '=#name?'
$wb.Names.Item("dbl_DBWeightP1").RefersTo = '=#name?'
$wb.Names.Item("dbl_DBWeightP2").RefersTo = '=#name?'
$wb.Names.Item("dbl_DBWeightP3").RefersTo = '=#name?'
$wb.Names.Item("dbl_DBWeightP4").RefersTo = '=#name?'
$wb.Names.Item("dbl_ScoreMax").RefersTo = '=#name?'
$wb.Names.Item("DynamicRange_CapBio").RefersTo = '=OFFSET(#name?, 0, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_CapGeo").RefersTo = '=OFFSET(#name?, 1, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_CapOther").RefersTo = '=OFFSET(#name?, 2, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_CapSmHy").RefersTo = '=OFFSET(#name?, 3, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_CapSolar").RefersTo = '=OFFSET(#name?, 4, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_CapWind").RefersTo = '=OFFSET(#name?, 5, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_GenBio").RefersTo = '=OFFSET(#name?, 0, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_GenGeo").RefersTo = '=OFFSET(#name?, 1, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_GenOther").RefersTo = '=OFFSET(#name?, 2, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_GenSmHy").RefersTo = '=OFFSET(#name?, 3, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_GenSolar").RefersTo = '=OFFSET(#name?, 4, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_GenWind").RefersTo = '=OFFSET(#name?, 5, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_InvBiofuels").RefersTo = '=OFFSET(#name?, 0, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_InvBiomass").RefersTo = '=OFFSET(#name?, 1, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_InvGeo").RefersTo = '=OFFSET(#name?, 2, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_InvOtherCE").RefersTo = '=OFFSET(#name?, 6, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_InvSmallDist").RefersTo = '=OFFSET(#name?, 7, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_InvSmHy").RefersTo = '=OFFSET(#name?, 3, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_InvSolar").RefersTo = '=OFFSET(#name?, 4, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_InvWind").RefersTo = '=OFFSET(#name?, 5, 1, 1, #name?)'
$wb.Names.Item("DynamicRange_Scores").RefersTo = '=OFFSET(#name?, 0, 0, 27, COUNTIF(#name?, "<>0"))'
$wb.Names.Item("int_Year").RefersTo = '=#name?'
$wb.Names.Item("list_Countries").RefersTo = '=#name?'
$wb.Names.Item("list_RegionOnGrid").RefersTo = '=#name?'
$wb.Names.Item("str_Country").RefersTo = '=#name?'
$wb.Names.Item("str_GridStatus").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink101").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink102").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink103").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink105").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink107").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink109").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink201").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink209").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink301").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink302").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink303").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink304").RefersTo = '=#name?'
$wb.Names.Item("str_Hyperlink305").RefersTo = '=#name?'
$wb.Names.Item("str_OptionInvType").RefersTo = '=#name?'
$wb.Names.Item("str_WeightScenario").RefersTo = '=#name?'

# --- Fix "Sri Lanka " (trailing space) -> "Sri Lanka" everywhere it appears ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A71").Value = "Sri Lanka"
$ws = $wb.Worksheets.Item(2)
$ws.Range("B59").Value = "Sri Lanka"
$ws = $wb.Worksheets.Item(3)
$ws.Range("B26").Value = "Sri Lanka"
$ws.Range("B88").Value = "Sri Lanka"
$ws.Range("B227").Value = "Sri Lanka"
$ws.Range("B243").Value = "Sri Lanka"
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2705").Value = "Sri Lanka"
$ws.Range("B2706").Value = "Sri Lanka"
$ws.Range("B2707").Value = "Sri Lanka"
$ws.Range("B2708").Value = "Sri Lanka"
$ws.Range("B2709").Value = "Sri Lanka"
$ws.Range("B2710").Value = "Sri Lanka"
$ws.Range("B2711").Value = "Sri Lanka"
$ws.Range("B2712").Value = "Sri Lanka"
$ws.Range("B2713").Value = "Sri Lanka"
$ws.Range("B2714").Value = "Sri Lanka"
$ws.Range("B2715").Value = "Sri Lanka"
$ws.Range("B2716").Value = "Sri Lanka"
$ws.Range("B2717").Value = "Sri Lanka"
$ws.Range("B2718").Value = "Sri Lanka"
$ws.Range("B2719").Value = "Sri Lanka"
$ws.Range("B2720").Value = "Sri Lanka"
$ws.Range("B2721").Value = "Sri Lanka"
$ws.Range("B2722").Value = "Sri Lanka"
$ws.Range("B2723").Value = "Sri Lanka"
$ws.Range("B2724").Value = "Sri Lanka"
$ws.Range("B2725").Value = "Sri Lanka"
$ws.Range("B2726").Value = "Sri Lanka"
$ws.Range("B2727").Value = "Sri Lanka"
$ws.Range("B2728").Value = "Sri Lanka"
$ws.Range("B2729").Value = "Sri Lanka"
$ws.Range("B2730").Value = "Sri Lanka"
$ws.Range("B2731").Value = "Sri Lanka"
$ws.Range("B2732").Value = "Sri Lanka"
$ws.Range("B2733").Value = "Sri Lanka"
$ws.Range("B2734").Value = "Sri Lanka"
$ws.Range("B2735").Value = "Sri Lanka"
$ws.Range("B2736").Value = "Sri Lanka"
$ws.Range("B2737").Value = "Sri Lanka"
$ws.Range("B2738").Value = "Sri Lanka"
$ws.Range("B2739").Value = "Sri Lanka"
$ws.Range("B2740").Value = "Sri Lanka"
$ws.Range("B2741").Value = "Sri Lanka"
$ws.Range("B2742").Value = "Sri Lanka"
$ws.Range("B2743").Value = "Sri Lanka"
$ws.Range("B2744").Value = "Sri Lanka"
$ws.Range("B2745").Value = "Sri Lanka"
$ws.Range("B2746").Value = "Sri Lanka"
$ws.Range("B2747").Value = "Sri Lanka"
$ws.Range("B2748").Value = "Sri Lanka"
$ws.Range("B2749").Value = "Sri Lanka"
$ws.Range("B2750").Value = "Sri Lanka"
$ws.Range("B2751").Value = "Sri Lanka"
$ws.Range("B2752").Value = "Sri Lanka"
$ws.Range("B2753").Value = "Sri Lanka"
$ws.Range("B2754").Value = "Sri Lanka"
$ws.Range("B2755").Value = "Sri Lanka"
$ws.Range("B2756").Value = "Sri Lanka"
$ws.Range("B2757").Value = "Sri Lanka"

# --- Tab colors: ARGB alpha byte 00 -> FF (opaque white) ---
foreach ($ws in $wb.Worksheets) {
    $ws.Tab.Color = 16777215
}

# --- Cell style rename: "Normal 2" -> "Excel Built-in Normal 2" ---
foreach ($s in $wb.Styles) {
    if ($s.Name -eq "Normal 2") {
        $s.Name = "Excel Built-in Normal 2"
    }
}

# --- Tab ratio: 427 -> 302 (window split between sheet tabs / scrollbar) ---
$wb.Worksheets.Item(1).Activate()
$excel.ActiveWindow.TabRatio = 0.302
